$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134; this shifts existing rows 134-153 down to 135-154
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data point
$ws.Range("A134").Value = 5
$ws.Range("B134").Value = "Macroferia Regional de Talca"
$ws.Range("C134").Value = "Maule"
$ws.Range("D134").Value = 44491
$ws.Range("E134").Value = 7
$ws.Range("F134").Value = 100112008
$ws.Range("G134").Value = "Coliflor"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 700
$ws.Range("L134").Value = 700
$ws.Range("M134").Value = 700
$ws.Range("N134").Value = "`$/unidad"
$ws.Range("O134").Value = "Región del Maule"
$ws.Range("P134").Value = 700
$ws.Range("Q134").Value = 1
$ws.Range("R134").Value = "Hortaliza"
